# Select the "HomePage" worksheet (sheet3.xml) which holds the locators
# table and is the tab that was active in the edited workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")
$ws.Activate()

# Fill in the new locator rows (5-8) that were previously blank.
# Values are written in the same order the shared-string table in the
# target workbook accumulated them.
$ws.Range("C5").Value = "search_query_top"
$ws.Range("A5").Value = "txt_search_bar"
$ws.Range("A6").Value = "lbl_search_result"
$ws.Range("C6").Value = "//ul[@class='product_list grid row']//a[contains(text(),'searchText')]"
$ws.Range("A7").Value = "btn_search"
$ws.Range("C7").Value = "submit_search"
$ws.Range("C8").Value = "//ul[@class='product_list grid row']//a[contains(text(),'Faded Short Sleeve T-shirts')]"
$ws.Range("A8").Value = "lbl_search_result_locator"

$ws.Range("B5").Value = "ID"
$ws.Range("B6").Value = "XPath"
$ws.Range("B7").Value = "Name"
$ws.Range("B8").Value = "XPath"

# Match the saved selection state from the diff (active cell A8).
$ws.Range("A8").Select()
